# MYCE version employed for SDEWES
# Update Project Info, Upgrades Info, Yearly Costs Info and Yearly Energy Averages
# to reflect the new (5-year) multi-year capacity-expansion run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Project Info" — refresh NPC / Operation cost / LCOE results
# ---------------------------------------------------------------------------
$wsProject = $wb.Worksheets.Item("Project Info")
$wsProject.Range("B2").Value = 1461480.5867019631
$wsProject.Range("B3").Value = 1182705.79454426
$wsProject.Range("B4").Value = 0.4772966583158948

# ---------------------------------------------------------------------------
# Sheet "Upgrades Info" — only one upgrade remains (drop Upgrade 2 / Upgrade 3)
# ---------------------------------------------------------------------------
$wsUpgrades = $wb.Worksheets.Item("Upgrades Info")
$wsUpgrades.Range("C1:D1").EntireColumn.Delete()

$wsUpgrades.Range("B2").Value = 6.1750338886374809
$wsUpgrades.Range("B3").Value = 144.21478917984061
$wsUpgrades.Range("B4").Value = 45.897096546460141
$wsUpgrades.Range("B5").Value = 0
$wsUpgrades.Range("B6").Value = 2762.957163131955
$wsUpgrades.Range("B7").Value = 60584.632934451052
$wsUpgrades.Range("B8").Value = 215427.20206011989
$wsUpgrades.Range("B9").Value = 0
$wsUpgrades.Range("B10").Value = 278774.79215770302

# ---------------------------------------------------------------------------
# Sheet "Yearly Costs Info" — extend from 3 to 5 years
# ---------------------------------------------------------------------------
$wsCosts = $wb.Worksheets.Item("Yearly Costs Info")
$wsCosts.Range("A5:A6").EntireRow.Insert()
$wsCosts.Range("A4").Copy()
$wsCosts.Range("A5").PasteSpecial(-4122)
$wsCosts.Range("A4").Copy()
$wsCosts.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsCosts.Range("A5").Value = "Year 4"
$wsCosts.Range("A6").Value = "Year 5"

$wsCosts.Range("B2").Value = 55.259143262639093
$wsCosts.Range("C2").Value = 6058.4632934451056
$wsCosts.Range("D2").Value = 4308.5440412023991
$wsCosts.Range("E2").Value = 10422.26647791014
$wsCosts.Range("F2").Value = 174013.2147760638
$wsCosts.Range("G2").Value = 72.376755136928622
$wsCosts.Range("H2").Value = 0

$wsCosts.Range("B3").Value = 55.259143262639093
$wsCosts.Range("C3").Value = 6058.4632934451056
$wsCosts.Range("D3").Value = 4308.5440412023991
$wsCosts.Range("E3").Value = 10422.26647791014
$wsCosts.Range("F3").Value = 198224.2912546304
$wsCosts.Range("G3").Value = 42.342319963365092
$wsCosts.Range("H3").Value = 0

$wsCosts.Range("B4").Value = 55.259143262639093
$wsCosts.Range("C4").Value = 6058.4632934451056
$wsCosts.Range("D4").Value = 4308.5440412023991
$wsCosts.Range("E4").Value = 10422.26647791014
$wsCosts.Range("F4").Value = 223257.02587713659
$wsCosts.Range("G4").Value = 21.39342939190227
$wsCosts.Range("H4").Value = 0

$wsCosts.Range("B5").Value = 55.259143262639093
$wsCosts.Range("C5").Value = 6058.4632934451056
$wsCosts.Range("D5").Value = 4308.5440412023991
$wsCosts.Range("E5").Value = 10422.26647791014
$wsCosts.Range("F5").Value = 252371.9633319019
$wsCosts.Range("G5").Value = 2.3731693227008011
$wsCosts.Range("H5").Value = 0

$wsCosts.Range("B6").Value = 55.259143262639093
$wsCosts.Range("C6").Value = 6058.4632934451056
$wsCosts.Range("D6").Value = 4308.5440412023991
$wsCosts.Range("E6").Value = 10422.26647791014
$wsCosts.Range("F6").Value = 282583.36114769487
$wsCosts.Range("G6").Value = 6.1200935434860062
$wsCosts.Range("H6").Value = 0

# ---------------------------------------------------------------------------
# Sheet "Yearly Energy Averages" — extend from 3 to 5 years
# ---------------------------------------------------------------------------
$wsEnergy = $wb.Worksheets.Item("Yearly Energy Averages")
$wsEnergy.Range("A5:A6").EntireRow.Insert()
$wsEnergy.Range("A4").Copy()
$wsEnergy.Range("A5").PasteSpecial(-4122)
$wsEnergy.Range("A4").Copy()
$wsEnergy.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsEnergy.Range("A5").Value = "Year 4"
$wsEnergy.Range("A6").Value = "Year 5"

$wsEnergy.Range("B2").Value = 0.21338634807427681
$wsEnergy.Range("C2").Value = 0.0029659324378359291
$wsEnergy.Range("D2").Value = 0.0013201352094334971

$wsEnergy.Range("B3").Value = 0.19702507381117679
$wsEnergy.Range("C3").Value = 0.00093369280609981458
$wsEnergy.Range("D3").Value = 0.00071160003392778468

$wsEnergy.Range("B4").Value = 0.18281428969469249
$wsEnergy.Range("C4").Value = 0.00015904585932589559
$wsEnergy.Range("D4").Value = 0.00033333122898987429

$wsEnergy.Range("B5").Value = 0.17042783136176221
$wsEnergy.Range("C5").Value = 0.0000003078223829393755
$wsEnergy.Range("D5").Value = 0.000034464480218878667

$wsEnergy.Range("B6").Value = 0.15958705292452219
$wsEnergy.Range("C6").Value = 0
$wsEnergy.Range("D6").Value = 0.000086741253239078702

# ---------------------------------------------------------------------------
# Cosmetic / locale cleanup picked up when the workbook was re-saved from an
# English Excel install (was previously saved from an Italian install):
# the built-in "Normal" cell style and the default Office theme name.
# ---------------------------------------------------------------------------
$wb.Styles.Item(1).Name = "Normal"
$wb.Theme.Name = "Office Theme"

Write-Output "done"
